$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A trial is now repeated if a stimulator error occurs -> update the
# generated stimulation order values in columns D (channels) and E (electrodes)
$ws.Range("E2").Value = "[(0,0)]"
$ws.Range("D2").Value = "[1, 2, 3]"

$ws.Range("E3").Value = "[(0,0)]"
$ws.Range("D3").Value = "[1]"

$ws.Range("E4").Value = "[(0,0)]"
$ws.Range("D4").Value = "[3]"

$ws.Range("E5").Value = "[(0,0)]"
$ws.Range("D5").Value = "[1, 2]"

# Move the active selection to D7 (was D12)
$ws.Range("D7").Select()
